$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.349.91'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.913.12'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.731'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '256.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.49%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.42'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.368'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '53.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0762'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0988'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.58%  '
$ws.Range('D14').Value = '2.190.11'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.736'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.12%  '
$ws.Range('E16').Value = '  +4.55%  '
$ws.Range('D17').Value = '1.934.39'
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').Value = '35.374.32'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '75.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.88%  '
$ws.Range('D20').Value = '0.0₃0849'
$ws.Range('E20').Value = '  +4.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '245.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('E22').Value = '  +6.21%  '
$ws.Range('E23').Value = '  +7.63%  '
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('E25').Value = '  +7.41%  '
$ws.Range('E26').Value = '  +3.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.58'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.85'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.21%  '
$ws.Range('E30').Value = '  +5.38%  '
$ws.Range('D31').Value = '4.128.99'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('E32').Value = '  +6.50%  '
$ws.Range('E33').Value = '  +24.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0593'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.02%  '
$ws.Range('E36').Value = '  +5.23%  '
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.919'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('E40').Value = '  +6.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.14'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.35%  '
$ws.Range('E43').Value = '  +3.55%  '
$ws.Range('E44').Value = '  +1.75%  '
$ws.Range('E45').Value = '  +5.40%  '
$ws.Range('D46').Value = '1.345.46'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('E48').Value = '  +4.50%  '
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.00%  '
$ws.Range('E51').Value = '  +6.75%  '
